$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Highlight the "Make a preview sprite for which tile/object you are placing"
# Beta Feedback entry (row 7) using the "Good" cell style.
$ws.Range("E7:F7").Style = "Good"

# ObjectMaps are no longer placeable in worldmaps, so the feedback item
# "Exiting level editor should take you to world menu" no longer applies.
# Remove it from the Beta Feedback list and pull the remaining entries up
# (the list continues after the section break at row 14, into row 15).
$ws.Range("E11:F11").ClearContents()

$ws.Range("E12").Value = "Exiting level editor should take you to world menu"
$ws.Range("F12").Value = 1

$ws.Range("E13").Value = "Add ability to add custom music to level"
$ws.Range("F13").Value = 1

$ws.Range("E15").Value = "Make right clicking editor tux allow you to change Tux's powerup state"
$ws.Range("F15").Value = 1

$ws.Range("E11").Select()
